# Updated estimates for revised paper.
# The table gains a new (4th) results column (column E), and the
# coefficient estimates in columns B and C are revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: turn the old text headers into the new numeric ones,
#     and extend the header into the new column E. -----------------
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 0

# New column E header - copy D1's formatting (bold font, thin border,
# centered/top alignment) so the new column matches the existing header
# row style, then set its value.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = 0

# --- Row 2: "FE" ----------------------------------------------------
$ws.Range("A2").Value = "FE"
$ws.Range("B2").Value = -0.28
$ws.Range("C2").Value = 0.66
$ws.Range("D2").Value = 0.2

# --- Row 3: "FE+Disg" -------------------------------------------------
$ws.Range("A3").Value = "FE+Disg"
$ws.Range("B3").Value = -0.28
$ws.Range("C3").Value = 0.18
$ws.Range("D3").Value = 0.2

# --- Row 4: "FE+Disg+Var" ---------------------------------------------
$ws.Range("A4").Value = "FE+Disg+Var"
$ws.Range("B4").Value = -0.28
$ws.Range("C4").Value = 0.18
$ws.Range("D4").Value = 0.2
